# Insert two new rows at the top of the data block (row 49) and shift
# existing rows 49:142 down to 51:144, then populate the two new rows
# with the new "Paine" records dated 44526.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 49; this shifts rows 49:142 -> 51:144
$ws.Range("A49:R50").Insert(-4121)   # -4121 = xlShiftDown

# Populate new row 49
$ws.Range("A49").Value = 11
$ws.Range("B49").Value = "Vega Monumental Concepción"
$ws.Range("C49").Value = "Bíobío"
$ws.Range("D49").Value = 44526
$ws.Range("E49").Value = 8
$ws.Range("F49").Value = 100112045
$ws.Range("G49").Value = "Zapallo"
$ws.Range("H49").Value = "Paine"
$ws.Range("I49").Value = "1a (guarda)"
$ws.Range("J49").Value = 400
$ws.Range("K49").Value = 140
$ws.Range("L49").Value = 150
$ws.Range("M49").Value = 145
$ws.Range("N49").Value = "$/kilo (volumen en unidades)"
$ws.Range("O49").Value = "Región de O'Higgins"
$ws.Range("P49").Value = 145
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"

# Populate new row 50
$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44526
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = 100112045
$ws.Range("G50").Value = "Zapallo"
$ws.Range("H50").Value = "Paine"
$ws.Range("I50").Value = "2a (guarda)"
$ws.Range("J50").Value = 200
$ws.Range("K50").Value = 100
$ws.Range("L50").Value = 100
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = "$/kilo (volumen en unidades)"
$ws.Range("O50").Value = "Región de O'Higgins"
$ws.Range("P50").Value = 100
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"

# Make sure date cells keep the date number format used elsewhere in column D
$ws.Range("D49:D50").NumberFormat = $ws.Range("D51").NumberFormat
